{"js": "// Replace the 25 \"A\u00d7B=C\" answer strings throughout the table with their\n// updated values, as described by the diff. Each old value is unique in\n// the document, so a simple body-wide search-and-replace per pair is safe.\nconst replacements = [\n  [\"114\u00d77=798\", \"437\u00d77=3059\"],\n  [\"563\u00d78=4504\", \"906\u00d77=6342\"],\n  [\"727\u00d78=5816\", \"599\u00d73=1797\"],\n  [\"664\u00d74=2656\", \"768\u00d77=5376\"],\n  [\"756\u00d79=6804\", \"262\u00d72=524\"],\n  [\"911\u00d75=4555\", \"977\u00d74=3908\"],\n  [\"886\u00d74=3544\", \"557\u00d79=5013\"],\n  [\"963\u00d72=1926\", \"604\u00d77=4228\"],\n  [\"493\u00d74=1972\", \"173\u00d76=1038\"],\n  [\"184\u00d73=552\", \"978\u00d73=2934\"],\n  [\"394\u00d75=1970\", \"404\u00d77=2828\"],\n  [\"822\u00d79=7398\", \"531\u00d76=3186\"],\n  [\"445\u00d78=3560\", \"290\u00d74=1160\"],\n  [\"976\u00d74=3904\", \"882\u00d72=1764\"],\n  [\"159\u00d77=1113\", \"974\u00d75=4870\"],\n  [\"251\u00d76=1506\", \"207\u00d76=1242\"],\n  [\"465\u00d72=930\", \"203\u00d79=1827\"],\n  [\"816\u00d73=2448\", \"201\u00d73=603\"],\n  [\"313\u00d72=626\", \"970\u00d75=4850\"],\n  [\"294\u00d78=2352\", \"255\u00d74=1020\"],\n  [\"146\u00d77=1022\", \"522\u00d78=4176\"],\n  [\"622\u00d73=1866\", \"759\u00d76=4554\"],\n  [\"794\u00d78=6352\", \"999\u00d76=5994\"],\n  [\"959\u00d79=8631\", \"946\u00d72=1892\"],\n  [\"615\u00d78=4920\", \"663\u00d79=5967\"],\n];\n\nconst body = context.document.body;\n\nfor (const [oldText, newText] of replacements) {\n  const found = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  found.load(\"items\");\n  await context.sync();\n\n  for (let i = 0; i < found.items.length; i++) {\n    found.items[i].insertText(newText, \"Replace\");\n  }\n  await context.sync();\n}\n", "ps1": "# Replace the 25 \"A\u00d7B=C\" answer strings throughout the table with their\n# updated values, as described by the diff. Each old value is unique in\n# the document, so a straightforward Find/Replace per pair is safe.\n$d = $word.ActiveDocument\n\n$replacements = @(\n    @(\"114\u00d77=798\", \"437\u00d77=3059\"),\n    @(\"563\u00d78=4504\", \"906\u00d77=6342\"),\n    @(\"727\u00d78=5816\", \"599\u00d73=1797\"),\n    @(\"664\u00d74=2656\", \"768\u00d77=5376\"),\n    @(\"756\u00d79=6804\", \"262\u00d72=524\"),\n    @(\"911\u00d75=4555\", \"977\u00d74=3908\"),\n    @(\"886\u00d74=3544\", \"557\u00d79=5013\"),\n    @(\"963\u00d72=1926\", \"604\u00d77=4228\"),\n    @(\"493\u00d74=1972\", \"173\u00d76=1038\"),\n    @(\"184\u00d73=552\", \"978\u00d73=2934\"),\n    @(\"394\u00d75=1970\", \"404\u00d77=2828\"),\n    @(\"822\u00d79=7398\", \"531\u00d76=3186\"),\n    @(\"445\u00d78=3560\", \"290\u00d74=1160\"),\n    @(\"976\u00d74=3904\", \"882\u00d72=1764\"),\n    @(\"159\u00d77=1113\", \"974\u00d75=4870\"),\n    @(\"251\u00d76=1506\", \"207\u00d76=1242\"),\n    @(\"465\u00d72=930\", \"203\u00d79=1827\"),\n    @(\"816\u00d73=2448\", \"201\u00d73=603\"),\n    @(\"313\u00d72=626\", \"970\u00d75=4850\"),\n    @(\"294\u00d78=2352\", \"255\u00d74=1020\"),\n    @(\"146\u00d77=1022\", \"522\u00d78=4176\"),\n    @(\"622\u00d73=1866\", \"759\u00d76=4554\"),\n    @(\"794\u00d78=6352\", \"999\u00d76=5994\"),\n    @(\"959\u00d79=8631\", \"946\u00d72=1892\"),\n    @(\"615\u00d78=4920\", \"663\u00d79=5967\")\n)\n\nforeach ($pair in $replacements) {\n    $oldText = $pair[0]\n    $newText = $pair[1]\n\n    $rng = $d.Content\n    $find = $rng.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Text = $oldText\n    $find.Replacement.Text = $newText\n    $find.Forward = $true\n    $find.Wrap = 1\n    $find.MatchCase = $true\n    $find.MatchWholeWord = $false\n    $find.MatchWildcards = $false\n    $find.Execute([ref]$oldText, [ref]$true, [ref]$false, [ref]$false, [ref]$false, [ref]$false, [ref]$true, [ref]1, [ref]$false, [ref]$newText, [ref]2)\n}\n"}
